# Add a new planned place (row 18) to the "Planned" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planned")

$ws.Cells.Item(18, 1).Value = "Iguazu Falls, Argentina"
$ws.Cells.Item(18, 2).Value = "[-25.69253235, -54.44111443902037]"
$ws.Cells.Item(18, 3).Value = "https://dynamic-media-cdn.tripadvisor.com/media/photo-o/2e/ed/a5/17/foz-do-iguacu.jpg"
